$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fila 5 - Camacho Peña,Carlos Andres
$ws.Range("B5").Value = 10
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = 10

# Fila 6 - Castro Obando,Sebastian
$ws.Range("B6").Value = 10
$ws.Range("C6").Value = 7
$ws.Range("D6").Value = 10

# Fila 7 - Chapid Tobar,Willian David
$ws.Range("B7").Value = 10
$ws.Range("C7").Value = 10
$ws.Range("D7").Value = 10

# Fila 9 - Gaitan Zambrano,Luis Felipe
$ws.Range("B9").Value = 9.4
$ws.Range("C9").Value = 10
$ws.Range("D9").Value = 9.4

# Fila 10 - Garcia Lopez,Jose Manuel
$ws.Range("B10").Value = 9.4
$ws.Range("D10").Value = 9.4

# Fila 12 - Gomez Valencia,Sebastian: se quita la observación
$ws.Range("E12").Value = ""

# Selección final tal como quedó en el archivo guardado
$ws.Range("E15").Select() | Out-Null
